$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column C (Author) for rows 2-30 with "Unknown" (no prior per-cell
# formatting existed there, so re-apply the default "Normal" cell style
# after writing the value so no stray style index gets stamped on the cell)
for ($r = 2; $r -le 30; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = "Unknown"
    $cell.Style = "Parasts"
}

# Update the view/selection to match the saved state
$ws.Activate()
$ws.Range("C2:C30").Select()
$excel.ActiveWindow.ScrollRow = 13
